$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append new row 26 to the Logs sheet
$ws.Range("A26").Value = "Afmelding nieuwsbrief"
$ws.Range("B26").Value = "mailmind.test@zohomail.eu"
$ws.Range("C26").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Range("D26").Value = "Afmelding"
$ws.Range("F26").Value = "2025-06-19 17:52:20"
$ws.Range("G26").Value = "Nee"

# Update the Dashboard summary count for "Afmelding"
$dash.Range("B3").Value = 7

# Expand the conditional formatting ranges to include the new row
$dRange = $ws.Range("D2:D25")
for ($i = 1; $i -le $dRange.FormatConditions.Count; $i++) {
    $dRange.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("D2:D26"))
}

$gRange = $ws.Range("G2:G25")
for ($i = 1; $i -le $gRange.FormatConditions.Count; $i++) {
    $gRange.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("G2:G26"))
}
